# Append the "nome app" / NOME APP section to the end of the document,
# right after the existing "Analisi della concorrenza" paragraph and
# before the section properties, exactly as produced by the target edit.
#
# We build the new paragraphs as a raw WordprocessingML fragment and
# insert it with Range.InsertXML at the very end of the document's
# story, which is the most reliable way to reproduce the exact run /
# paragraph-mark structure (including an empty spacer paragraph that
# has no run at all, a proofErr-wrapped "del team" span, a
# yellow-highlighted "nome app" / "NOME APP" run, and the centered
# heading paragraph) shown in the diff.

$d = $word.ActiveDocument

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body><w:p>
  <w:pPr>
    <w:rPr>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
    <w:t xml:space="preserve">Considerati i bisogni degli utenti, le competenze </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
    <w:t>del team</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
    <w:t xml:space="preserve"> e le opportunità offerte dal mercato, emerge la necessità di un sistema che renda l’accesso ai servizi sanitari più intuitivo ed efficiente. Da qui la scelta di </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
    <w:t>“</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
      <w:highlight w:val="yellow"/>
    </w:rPr>
    <w:t>nome app</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
    <w:t>”</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
    <w:t>, un nome che riflette l’idea di una sanità digitale intelligente, semplice e accessibile a tutti.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="center"/>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="36"/>
      <w:szCs w:val="36"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="36"/>
      <w:szCs w:val="36"/>
      <w:highlight w:val="yellow"/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t>NOME APP</w:t>
  </w:r>
</w:p></w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$target = $d.Range($d.Content.End, $d.Content.End)
$target.InsertXML($xml)
